$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Added Minifit Jr. to BOM: new Digi-Key part-number/link pair on row 2
# (ACC*, DUE, SERVO* / Molex 5566-2 Mini FIT connector row)
$ws.Range("F2").Value = "WM3875-ND"
$ws.Range("G2").Value = "http://www.digikey.com/product-detail/en/0039299027/WM3875-ND/2002681"

# Hyperlink the new Digi-Key link cell
$ws.Hyperlinks.Add($ws.Range("G2"), "http://www.digikey.com/product-detail/en/0039299027/WM3875-ND/2002681")
$g2Font = $ws.Range("G2").Font
$g2Font.Underline = $false
$g2Font.Color = 16711680

# Trim stray trailing spaces from existing Digi-Key part numbers
$ws.Range("F5").Value = "399-8269-1-ND"
$ws.Range("F7").Value = "SMBJ8.5ALFCT-ND"
$ws.Range("F8").Value = "1SMB5913BT3GOSCT-ND"
$ws.Range("F9").Value = "WK6265-ND"
$ws.Range("F10").Value = "160-1169-1-ND"
$ws.Range("F11").Value = "RMCF1206JT1K00CT-ND"
$ws.Range("F12").Value = "RMCF1206FT15K4CT-ND"

# Restore selection cursor to F3
$ws.Range("F3").Select()
